$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the player on row 13 (Keita -> Fabinho) and update his stats
$ws.Range("A13").Value = "Fabinho"
$ws.Range("B13").Value = 2
$ws.Range("D13").Value = 1989
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 75

# Move the selection to F14 (simulating the user tabbing/entering down the row)
# and scroll the view back to the top-left (A1).
[void]$ws.Range("A1").Select()
[void]$ws.Range("F14").Select()
